# Update price (D) and volume-change (E) columns per latest crypto data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.423.16"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.673.16"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.32"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5342"

$ws.Range("E7").Value = "  +0.49%  "

$ws.Range("E8").Value = "  +1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06388"
$ws.Range("E9").Value = "  +1.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.88"
$ws.Range("E10").Value = "  +2.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07867"
$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.535"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.675.01"
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.902.88"
$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5618"
$ws.Range("E15").Value = "  +2.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8198"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.23"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.441.13"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.729"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "197.93"
$ws.Range("E21").Value = "  +3.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.31"
$ws.Range("E22").Value = "  +2.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.082"
$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("E24").Value = "  +0.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.24"
$ws.Range("E25").Value = "  +0.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1228"
$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.262"
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.27"
$ws.Range("E28").Value = "  +1.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").Value = "  +2.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05930"
$ws.Range("E30").Value = "  +3.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  +1.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.562"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.337"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.613"
$ws.Range("E34").Value = "  +1.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9696"
$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5838"
$ws.Range("E38").Value = "  +2.00%  "

$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.080.15"
$ws.Range("E40").Value = "  +4.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.932"
$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8652"
$ws.Range("E42").Value = "  +1.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.21"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.62"
$ws.Range("E46").Value = "  +3.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.014"
$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4417"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.036"
$ws.Range("E50").Value = "  +2.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05163"
$ws.Range("E51").Value = "  +0.16%  "
